$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-4 down to 4-5
$ws.Rows.Item(3).Insert()

# Fill in the new row 3 with the inserted record (date 2022-05-27 = 44708)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44708
$ws.Range("D3").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100108
$ws.Range("H3").Value = "Tropicales y subtropicales"
$ws.Range("I3").Value = 100108007
$ws.Range("J3").Value = "Coco"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 21000
$ws.Range("P3").Value = 20500
$ws.Range("Q3").Value = "$/malla 20 unidades"
$ws.Range("R3").Value = "Perú"
$ws.Range("S3").Value = 1025
$ws.Range("T3").Value = 20
